# Update "want-to-go" counts (column F) across the workbook's sheets
# to reflect the latest scrape, per commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws4 = $wb.Worksheets.Item("全部类型")

# --- Sheet: 展览 (exhibitions) ---
$ws1.Range("F2").Value  = 218
$ws1.Range("F3").Value  = 121
$ws1.Range("F5").Value  = 1004
$ws1.Range("F6").Value  = 5605
$ws1.Range("F7").Value  = 506
$ws1.Range("F8").Value  = 716
$ws1.Range("F17").Value = 1887
$ws1.Range("F18").Value = 1483
$ws1.Range("F19").Value = 949
$ws1.Range("F20").Value = 299
$ws1.Range("F22").Value = 345
$ws1.Range("F23").Value = 565
$ws1.Range("F24").Value = 164
$ws1.Range("F25").Value = 1059
$ws1.Range("F28").Value = 3037
$ws1.Range("F29").Value = 183
$ws1.Range("F33").Value = 41
$ws1.Range("F34").Value = 417
$ws1.Range("F39").Value = 301
$ws1.Range("F40").Value = 746
$ws1.Range("F42").Value = 57

# --- Sheet: 演出 (performances) ---
$ws2.Range("F2").Value = 34
$ws2.Range("F4").Value = 208

# --- Sheet: 全部类型 (all types) ---
$ws4.Range("F3").Value  = 218
$ws4.Range("F4").Value  = 121
$ws4.Range("F5").Value  = 1004
$ws4.Range("F6").Value  = 34
$ws4.Range("F7").Value  = 5605
$ws4.Range("F8").Value  = 506
$ws4.Range("F9").Value  = 716
$ws4.Range("F11").Value = 208
$ws4.Range("F23").Value = 1887
$ws4.Range("F24").Value = 1483
$ws4.Range("F25").Value = 949
$ws4.Range("F27").Value = 345
$ws4.Range("F29").Value = 565
$ws4.Range("F30").Value = 164
$ws4.Range("F31").Value = 1059
$ws4.Range("F32").Value = 3037
$ws4.Range("F33").Value = 183
$ws4.Range("F37").Value = 41
$ws4.Range("F38").Value = 417
$ws4.Range("F42").Value = 301
$ws4.Range("F43").Value = 746
